$p = $ppt.ActivePresentation

# --- Update existing slides 1-10 (word swap + page range 65-66 -> 63-64) ---
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "豊か"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "ゆたか"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "abundant, plentiful, rich, ample | rich, wealthy, affluent, well-off | open (mind), relaxed, easy | plump (e.g. breasts),..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

$s = $p.Slides.Item(2)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "豊作"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "ほうさく"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "abundant harvest, bumper crop..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

$s = $p.Slides.Item(3)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "得る"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "える"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "to get, to earn, to acquire, to procure, to gain, to secure, to attain, to obtain, to win | to understand, to comprehend ..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "得意"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "とくい"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "satisfaction, pride, triumph, elation | one's strong point, one's forte, one's specialty | regular customer, regular clie..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "逆らう"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "さからう"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "to go against, to oppose, to disobey, to defy..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

$s = $p.Slides.Item(6)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "逆"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "ぎゃく"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "reverse, opposite | converse (of a hypothesis, etc.) | inverse (function)..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

$s = $p.Slides.Item(7)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "お互い"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "おたがい"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "each other, one another..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

$s = $p.Slides.Item(8)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "相互"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "そうご"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "mutual, reciprocal..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

$s = $p.Slides.Item(9)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "伸びる"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "のびる"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "to stretch, to extend, to lengthen, to grow (of hair, height, grass, etc.) | to straighten out, to be flattened, to becom..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

$s = $p.Slides.Item(10)
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "与える"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "あたえる"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "to give (esp. to someone of lower status), to bestow, to grant, to confer, to present, to award | to provide, to afford, ..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

# --- Add 3 new slides (11-13), duplicated from slide 10 layout/format ---
$src = $p.Slides.Item(10)
$src.Duplicate() | Out-Null
$p.Slides.Item($p.Slides.Count).Duplicate() | Out-Null
$p.Slides.Item($p.Slides.Count).Duplicate() | Out-Null

$s = $p.Slides.Item(11)
$s.Name = "Slide 11"
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "貸与"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "たいよ"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "loan, lending..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

$s = $p.Slides.Item(12)
$s.Name = "Slide 12"
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "可能"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "かのう"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "possible, potential, practicable, feasible..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"

$s = $p.Slides.Item(13)
$s.Name = "Slide 13"
$s.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "趣味"
$s.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "しゅみ"
$s.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = "hobby, pastime | tastes, preference, liking..."
$s.Shapes.Item(4).TextFrame.TextRange.Runs(1).Text = "63-64"
